$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new rows of data (činnost / hodiny) below the existing entries.
# Values are entered in this particular order so the shared-strings table
# ends up with the same ordering as the original author's edit.
$ws.Range("A15").Value = "Úpravy implementací"
$ws.Range("B15").Value = 2

$ws.Range("A17").Value = "Dodělání 4.iterace"
$ws.Range("B17").Value = 0.5

$ws.Range("A16").Value = "Testování, napsání zprávy"
$ws.Range("B16").Value = 2

# Update the total formula to reflect the extended range (matches existing pattern shift)
$ws.Range("B4").Formula = "=SUM(B7:B1001)"

# Move the active selection to E16 (matches final saved cursor position)
$ws.Range("E16").Select()

$wb.Save()
